$d = $word.ActiveDocument

# --- Change 1: merge the two runs of "Make sure we our happy ... highest
# mark" / " we can" into a single run, dropping the stray _GoBack bookmark
# that sat between them.
$d.Content.Find.Execute(
    "Make sure we our happy with the quality of each deliverable to achieve the highest mark we can",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Make sure we our happy with the quality of each deliverable to achieve the highest mark we can",
    2) | Out-Null

# --- Change 2: the "Clear Stage 4 backlog" row of the Gantt table is done,
# remove it entirely. Word then re-anchors the (single, reserved) _GoBack
# bookmark at the start of the next row's first cell ("Create HR database").
$t = $d.Tables.Item(1)
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "Clear Stage 4 backlog*") {
        $row.Delete()
        break
    }
}

for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "Create HR database*") {
        $cellRange = $row.Cells.Item(1).Range
        $startPoint = $d.Range($cellRange.Start, $cellRange.Start)
        $d.Bookmarks.Add("_GoBack", $startPoint) | Out-Null
        break
    }
}
